# Weekly CompStat (6th Precinct) refresh: new report week + updated crime counts.
# Maps each changed cell from the source workbook diff onto COM Range/Cells writes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Report header
$ws.Cells.Item(8,1).Value = "Volume 31   Number  5"  # A8

# Row 9: Report date range
$ws.Cells.Item(9,3).Value = "Report Covering the Week  1/29/2024  Through  2/4/2024"  # C9

# Row 15: 'Rape' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(15,7).Value = "'0"  # G15
$ws.Cells.Item(15,7).NumberFormat = "General"
$ws.Cells.Item(15,8).Value = "***.*"  # H15
$ws.Cells.Item(15,8).NumberFormat = "General"
$ws.Cells.Item(15,13).Value = -100  # M15
$ws.Cells.Item(15,13).NumberFormat = "#,##0.0;`"-`"#,##0.0"

# Row 16: 'Robbery' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(16,3).Value = 1  # C16
$ws.Cells.Item(16,4).Value = 3  # D16
$ws.Cells.Item(16,5).Value = -66.666666666666  # E16
$ws.Cells.Item(16,6).Value = 10  # F16
$ws.Cells.Item(16,7).Value = 11  # G16
$ws.Cells.Item(16,8).Value = -9.090909090909  # H16
$ws.Cells.Item(16,9).Value = 12  # I16
$ws.Cells.Item(16,10).Value = 17  # J16
$ws.Cells.Item(16,11).Value = -29.411764705882  # K16
$ws.Cells.Item(16,12).Value = -40  # L16
$ws.Cells.Item(16,13).Value = -40  # M16
$ws.Cells.Item(16,14).Value = -86.206896551724  # N16

# Row 17: 'Fel. Assault' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(17,3).Value = "'0"  # C17
$ws.Cells.Item(17,3).NumberFormat = "General"
$ws.Cells.Item(17,4).Value = 3  # D17
$ws.Cells.Item(17,5).Value = -100  # E17
$ws.Cells.Item(17,6).Value = 6  # F17
$ws.Cells.Item(17,7).Value = 8  # G17
$ws.Cells.Item(17,8).Value = -25  # H17
$ws.Cells.Item(17,10).Value = 12  # J17
$ws.Cells.Item(17,11).Value = -33.333333333333  # K17
$ws.Cells.Item(17,12).Value = 0  # L17
$ws.Cells.Item(17,13).Value = 14.285714285714  # M17
$ws.Cells.Item(17,14).Value = -57.894736842105  # N17

# Row 18: 'Burglary' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(18,3).Value = 3  # C18
$ws.Cells.Item(18,4).Value = 4  # D18
$ws.Cells.Item(18,5).Value = -25  # E18
$ws.Cells.Item(18,6).Value = 21  # F18
$ws.Cells.Item(18,7).Value = 26  # G18
$ws.Cells.Item(18,8).Value = -19.230769230769  # H18
$ws.Cells.Item(18,9).Value = 23  # I18
$ws.Cells.Item(18,10).Value = 35  # J18
$ws.Cells.Item(18,11).Value = -34.285714285714  # K18
$ws.Cells.Item(18,12).Value = -14.814814814814  # L18
$ws.Cells.Item(18,13).Value = -17.857142857142  # M18
$ws.Cells.Item(18,14).Value = -67.605633802816  # N18

# Row 19: 'Gr. Larceny' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(19,3).Value = 18  # C19
$ws.Cells.Item(19,5).Value = -25  # E19
$ws.Cells.Item(19,7).Value = 106  # G19
$ws.Cells.Item(19,8).Value = -26.415094339622  # H19
$ws.Cells.Item(19,9).Value = 95  # I19
$ws.Cells.Item(19,10).Value = 132  # J19
$ws.Cells.Item(19,11).Value = -28.030303030303  # K19
$ws.Cells.Item(19,12).Value = 2.150537634408  # L19
$ws.Cells.Item(19,13).Value = 6.741573033707  # M19
$ws.Cells.Item(19,14).Value = -58.333333333333  # N19

# Row 20: 'G.L.A.' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(20,3).Value = "'0"  # C20
$ws.Cells.Item(20,3).NumberFormat = "General"
$ws.Cells.Item(20,5).Value = -100  # E20
$ws.Cells.Item(20,6).Value = 1  # F20
$ws.Cells.Item(20,8).Value = -50  # H20
$ws.Cells.Item(20,10).Value = 3  # J20
$ws.Cells.Item(20,11).Value = 0  # K20
$ws.Cells.Item(20,13).Value = -25  # M20
$ws.Cells.Item(20,14).Value = -96.428571428571  # N20

# Row 21: 'TOTAL' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(21,3).Value = 22  # C21
$ws.Cells.Item(21,4).Value = 35  # D21
$ws.Cells.Item(21,5).Value = -37.142857142857  # E21
$ws.Cells.Item(21,6).Value = 116  # F21
$ws.Cells.Item(21,7).Value = 153  # G21
$ws.Cells.Item(21,8).Value = -24.183006535947  # H21
$ws.Cells.Item(21,9).Value = 141  # I21
$ws.Cells.Item(21,10).Value = 200  # J21
$ws.Cells.Item(21,11).Value = -29.5  # K21
$ws.Cells.Item(21,12).Value = -9.615384615384  # L21
$ws.Cells.Item(21,13).Value = -5.369127516778  # M21
$ws.Cells.Item(21,14).Value = -71.224489795918  # N21

# Row 22: 'Transit' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(22,4).Value = 2  # D22
$ws.Cells.Item(22,6).Value = 2  # F22
$ws.Cells.Item(22,7).Value = 5  # G22
$ws.Cells.Item(22,8).Value = -60  # H22
$ws.Cells.Item(22,10).Value = 5  # J22
$ws.Cells.Item(22,11).Value = 0  # K22

# Row 24: 'Petit Larceny' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(24,3).Value = 50  # C24
$ws.Cells.Item(24,4).Value = 37  # D24
$ws.Cells.Item(24,5).Value = 35.135135135135  # E24
$ws.Cells.Item(24,6).Value = 121  # F24
$ws.Cells.Item(24,7).Value = 150  # G24
$ws.Cells.Item(24,8).Value = -19.333333333333  # H24
$ws.Cells.Item(24,9).Value = 150  # I24
$ws.Cells.Item(24,10).Value = 183  # J24
$ws.Cells.Item(24,11).Value = -18.032786885245  # K24
$ws.Cells.Item(24,12).Value = -5.66037735849  # L24
$ws.Cells.Item(24,13).Value = 4.895104895104  # M24

# Row 25: 'Misd. Assault' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(25,3).Value = 8  # C25
$ws.Cells.Item(25,4).Value = 6  # D25
$ws.Cells.Item(25,5).Value = 33.333333333333  # E25
$ws.Cells.Item(25,7).Value = 27  # G25
$ws.Cells.Item(25,8).Value = 3.703703703703  # H25
$ws.Cells.Item(25,9).Value = 36  # I25
$ws.Cells.Item(25,10).Value = 37  # J25
$ws.Cells.Item(25,11).Value = -2.702702702702  # K25
$ws.Cells.Item(25,12).Value = 111.764705882353  # L25
$ws.Cells.Item(25,13).Value = 140  # M25

# Row 26: 'UCR Rape*' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(26,7).Value = "'0"  # G26
$ws.Cells.Item(26,7).NumberFormat = "General"
$ws.Cells.Item(26,8).Value = "***.*"  # H26
$ws.Cells.Item(26,8).NumberFormat = "General"

# Row 27: 'Other Sex Crimes' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(27,4).Value = 1  # D27
$ws.Cells.Item(27,7).Value = 7  # G27
$ws.Cells.Item(27,8).Value = -71.428571428571  # H27
$ws.Cells.Item(27,10).Value = 9  # J27
$ws.Cells.Item(27,11).Value = -77.777777777777  # K27

# Row 30: 'Hate Crimes' data refresh (week/28-day/YTD counts + %Chg)
$ws.Cells.Item(30,3).Value = 1  # C30
$ws.Cells.Item(30,3).NumberFormat = "#,##0"
$ws.Cells.Item(30,6).Value = 1  # F30
$ws.Cells.Item(30,6).NumberFormat = "#,##0"
$ws.Cells.Item(30,9).Value = 1  # I30
$ws.Cells.Item(30,9).NumberFormat = "#,##0"
